$p = $ppt.ActivePresentation
Write-Host "Designs before: $($p.Designs.Count)"
$m = $p.Slides.Item(1).Master
try {
  $m.Delete()
  Write-Host "deleted master"
} catch {
  Write-Host "delete err: $_"
}
Write-Host "Designs after: $($p.Designs.Count)"
